$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.808.13"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").Value = "1.637.88"
$ws.Range("E3").Value = "  +0.80%  "
$ws.Range("E4").Value = "  +0.42%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.37"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.40%  "
$ws.Range("E6").Value = "  -0.75%  "
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "28.83"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.260"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0608"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0898"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.65%  "
$ws.Range("D12").Value = "1.872.61"
$ws.Range("E12").Value = "  +0.88%  "
$ws.Range("D13").Value = "1.636.59"
$ws.Range("E13").Value = "  +1.14%  "
$ws.Range("E14").Value = "  +3.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "9.56"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +7.94%  "
$ws.Range("E16").Value = "  -0.49%  "
$ws.Range("D17").Value = "29.809.97"
$ws.Range("E17").Value = "  -0.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.24"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "237.82"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.75%  "
$ws.Range("D20").Value = "0.0₃0701"
$ws.Range("E20").Value = "  -0.69%  "
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.89"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.60%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.16"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.39"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.58"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.85%  "
$ws.Range("E27").Value = "  -1.53%  "
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("E29").Value = "  +0.39%  "
$ws.Range("E30").Value = "  +0.93%  "
$ws.Range("E31").Value = "  -0.83%  "
$ws.Range("E32").Value = "  +0.94%  "
$ws.Range("E33").Value = "  -1.32%  "
$ws.Range("D34").Value = "1.417.53"
$ws.Range("E34").Value = "  -0.60%  "
$ws.Range("E35").Value = "  +1.85%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.01"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.69"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -6.12%  "
$ws.Range("E38").Value = "  +1.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "76.48"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +10.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.565"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.00%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("E43").Value = "  -0.28%  "
$ws.Range("E44").Value = "  -2.94%  "
$ws.Range("E45").Value = "  +0.51%  "
$ws.Range("E46").Value = "  -1.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "49.93"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -8.91%  "
$ws.Range("D48").Value = "1.781.04"
$ws.Range("E48").Value = "  +0.91%  "
$ws.Range("E49").Value = "  -1.54%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "93.33"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +5.20%  "
$ws.Range("E51").Value = "  +0.61%  "
